# viz updates, three more langs added for llama
# Fill in accuracy (col B) / uncertainty (col C) results for the languages
# that previously only had placeholder "10.0" values and an empty
# uncertainty cell (rows 16-20 on the "llama" sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("llama")

$ws.Range("B16").Value = 82.72
$ws.Range("C16").Value = 0.56

$ws.Range("B17").Value = 84.01
$ws.Range("C17").Value = 0.52

$ws.Range("B18").Value = 84.06
$ws.Range("C18").Value = 0.49

$ws.Range("B19").Value = 75.74
$ws.Range("C19").Value = 1.21

$ws.Range("B20").Value = 60.01
$ws.Range("C20").Value = 1.28

# These uncertainty cells were previously blank (right-aligned style without
# the explicit readingOrder attribute); align them with the rest of the
# uncertainty column (right-aligned, same style as C2:C15) now that they
# hold real values.
$ws.Range("C16:C20").HorizontalAlignment = -4152
